$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new 2022-Q4 row at the top of the data
#    and shift the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Give the brand-new A8 cell the same look (border/bold/alignment) as the
# other index cells in column A before we put a value into it.
$summary.Range("A7").Copy()
$summary.Range("A8").PasteSpecial(-4122)

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 1.88

$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 2.1

$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 6
$summary.Range("D4").Value = 5.69

$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 7
$summary.Range("D5").Value = 2.46

$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 9
$summary.Range("D6").Value = 0.76

$summary.Range("B7").Value = "2021-Q1"
$summary.Range("C7").Value = 19
$summary.Range("D7").Value = 0.98

$summary.Range("A8").Value = 6
$summary.Range("B8").Value = "2020-Q4"
$summary.Range("C8").Value = 4
$summary.Range("D8").Value = 0.02

# ---------------------------------------------------------------------------
# 2) Add a new "2022-Q4" sheet right after "总计" and before "2022-Q3",
#    pushing every other quarter sheet one position to the right.
#    Duplicating the existing "2022-Q3" sheet gives us identical sheet
#    formatting (borders, fonts, page margins) for free.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $summary)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# The template sheet only has 2 data rows (rows 2-3); we need 7 (rows 2-8).
# Stamp column A's styled look onto the extra rows before filling them in.
$q4.Range("A2").Copy()
$q4.Range("A4:A8").PasteSpecial(-4122)

# Header row (already copied from the template, just keep the same text).
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Row 2
$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'011174"
$q4.Range("C2").Value = "中庚价值品质一年持有期混合"
$q4.Range("D2").Value = "'67.05"
$q4.Range("E2").Value = "'93.59"
$q4.Range("F2").Value = "'2.67"
$q4.Range("G2").Value = "'1.7902"
$q4.Range("H2").Value = 9

# Row 3
$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'257050"
$q4.Range("C3").Value = "国联安主题驱动混合"
$q4.Range("D3").Value = "'1.45"
$q4.Range("E3").Value = "'94.07"
$q4.Range("F3").Value = "'5.98"
$q4.Range("G3").Value = "'0.0867"
$q4.Range("H3").Value = 3

# Row 4
$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'850007"
$q4.Range("C4").Value = "海通智选一年持有期股票B"
$q4.Range("D4").Value = "'0.30"
$q4.Range("E4").Value = "'82.33"
$q4.Range("F4").Value = "'0.75"
$q4.Range("G4").Value = "'0.0022"
$q4.Range("H4").Value = 4

# Row 5
$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'850788"
$q4.Range("C5").Value = "海通智选一年持有期股票A"
$q4.Range("D5").Value = "'0.20"
$q4.Range("E5").Value = "'82.33"
$q4.Range("F5").Value = "'0.75"
$q4.Range("G5").Value = "'0.0015"
$q4.Range("H5").Value = 4

# Row 6
$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'519222"
$q4.Range("C6").Value = "海富通欣益灵活配置混合A"
$q4.Range("D6").Value = "'0.25"
$q4.Range("E6").Value = "'31.65"
$q4.Range("F6").Value = "'0.16"
$q4.Range("G6").Value = "'0.0004"
$q4.Range("H6").Value = 9

# Row 7
$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'519221"
$q4.Range("C7").Value = "海富通欣益灵活配置混合C"
$q4.Range("D7").Value = "'0.10"
$q4.Range("E7").Value = "'31.65"
$q4.Range("F7").Value = "'0.16"
$q4.Range("G7").Value = "'0.0002"
$q4.Range("H7").Value = 9

# Row 8
$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "'850799"
$q4.Range("C8").Value = "海通智选一年持有期股票C"
$q4.Range("D8").Value = "'0.00"
$q4.Range("E8").Value = "'82.33"
$q4.Range("F8").Value = "'0.75"
$q4.Range("G8").Value = 0
$q4.Range("H8").Value = 4
